# Initial Data File Update
# Adds two new transaction rows (35 and 36) to the "Transacciones" sheet,
# mirroring the rows already present for the "Extra" / "Golosinas" category,
# and updates the sheet's selection/view state to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transacciones")

# --- Row 35: Galletas Emperador Chocolate -------------------------------
$ws.Range("A35").Value = 43567
$ws.Range("A34").Copy()
$ws.Range("A35").PasteSpecial(-4122)   # xlPasteFormats - reuse the existing date style

$ws.Range("B35").Value = 15
$ws.Range("C35").Value = "Galletas Emperador Chocolate"
$ws.Range("D35").Value = "Golosinas"
$ws.Range("E35").Value = "Gasto"
$ws.Range("F35").Value = "Tarjeta Santander"
$ws.Range("G35").Value = "Extra"

$ws.Range("K35").Value = 7358.64
$ws.Range("L35").Formula = "=L34-B35"
$ws.Range("M35").Value = 504
$ws.Range("N35").Formula = "=SUM(K35:M35)"
$ws.Range("O35").Formula = "=N35-4000"

# --- Row 36: Licuado ------------------------------------------------------
$ws.Range("A36").Value = 43567
$ws.Range("A34").Copy()
$ws.Range("A36").PasteSpecial(-4122)   # xlPasteFormats - reuse the existing date style

$ws.Range("B36").Value = 13.5
$ws.Range("C36").Value = "Licuado"
$ws.Range("D36").Value = "Golosinas"
$ws.Range("E36").Value = "Gasto"
$ws.Range("F36").Value = "Tarjeta Santander"
$ws.Range("G36").Value = "Extra"

$ws.Range("K36").Value = 7358.64
$ws.Range("L36").Formula = "=L35-B36"
$ws.Range("M36").Value = 504
$ws.Range("N36").Formula = "=SUM(K36:M36)"
$ws.Range("O36").Formula = "=N36-4000"

# --- View state: move the active selection like the source file did -----
$ws.Range("P37").Select()
